$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, [string]$text) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value2 = $text
    $cell.Style = $origStyle
}

# Row 2
Set-TextValue $ws.Cells.Item(2, 4) "24.923.68"
Set-TextValue $ws.Cells.Item(2, 5) "  +2.17%  "

# Row 3
Set-TextValue $ws.Cells.Item(3, 4) "1.707.08"
Set-TextValue $ws.Cells.Item(3, 5) "  +1.74%  "

# Row 4
Set-TextValue $ws.Cells.Item(4, 5) "  -0.10%  "

# Row 5
Set-TextValue $ws.Cells.Item(5, 4) "312.88"
Set-TextValue $ws.Cells.Item(5, 5) "  +2.22%  "

# Row 6
Set-TextValue $ws.Cells.Item(6, 4) "0.9998"
Set-TextValue $ws.Cells.Item(6, 5) "  +0.02%  "

# Row 7
Set-TextValue $ws.Cells.Item(7, 4) "0.3742"
Set-TextValue $ws.Cells.Item(7, 5) "  +1.19%  "

# Row 8
Set-TextValue $ws.Cells.Item(8, 4) "49.41"
Set-TextValue $ws.Cells.Item(8, 5) "  +3.74%  "

# Row 9
Set-TextValue $ws.Cells.Item(9, 4) "0.3441"
Set-TextValue $ws.Cells.Item(9, 5) "  +0.29%  "

# Row 10
Set-TextValue $ws.Cells.Item(10, 4) "1.227"
Set-TextValue $ws.Cells.Item(10, 5) "  +5.48%  "

# Row 11
Set-TextValue $ws.Cells.Item(11, 4) "0.07540"
Set-TextValue $ws.Cells.Item(11, 5) "  +4.14%  "

# Row 12
Set-TextValue $ws.Cells.Item(12, 4) "1.002"
Set-TextValue $ws.Cells.Item(12, 5) "  -0.09%  "

# Row 13
Set-TextValue $ws.Cells.Item(13, 4) "21.18"
Set-TextValue $ws.Cells.Item(13, 5) "  +5.22%  "

# Row 14
Set-TextValue $ws.Cells.Item(14, 4) "6.332"
Set-TextValue $ws.Cells.Item(14, 5) "  +3.08%  "

# Row 15
Set-TextValue $ws.Cells.Item(15, 4) "7.065"
Set-TextValue $ws.Cells.Item(15, 5) "  +4.86%  "

# Row 16
Set-TextValue $ws.Cells.Item(16, 4) "1.706.55"
Set-TextValue $ws.Cells.Item(16, 5) "  +1.84%  "

# Row 17
Set-TextValue $ws.Cells.Item(17, 4) "0.00001131"
Set-TextValue $ws.Cells.Item(17, 5) "  +2.67%  "

# Row 18
Set-TextValue $ws.Cells.Item(18, 4) "0.06732"
Set-TextValue $ws.Cells.Item(18, 5) "  +0.97%  "

# Row 19
Set-TextValue $ws.Cells.Item(19, 4) "0.9990"
Set-TextValue $ws.Cells.Item(19, 5) "  -0.02%  "

# Row 20
Set-TextValue $ws.Cells.Item(20, 4) "83.99"
Set-TextValue $ws.Cells.Item(20, 5) "  +3.88%  "

# Row 21
Set-TextValue $ws.Cells.Item(21, 4) "17.36"
Set-TextValue $ws.Cells.Item(21, 5) "  +5.53%  "

# Row 22
Set-TextValue $ws.Cells.Item(22, 4) "6.386"
Set-TextValue $ws.Cells.Item(22, 5) "  +4.67%  "

# Row 23
Set-TextValue $ws.Cells.Item(23, 4) "13.13"
Set-TextValue $ws.Cells.Item(23, 5) "  +8.18%  "

# Row 24
Set-TextValue $ws.Cells.Item(24, 4) "24.918.33"
Set-TextValue $ws.Cells.Item(24, 5) "  +2.28%  "

# Row 25
Set-TextValue $ws.Cells.Item(25, 4) "2.449"
Set-TextValue $ws.Cells.Item(25, 5) "  +0.43%  "

# Row 26
Set-TextValue $ws.Cells.Item(26, 5) "  +5.57%  "

# Row 27
Set-TextValue $ws.Cells.Item(27, 4) "20.40"
Set-TextValue $ws.Cells.Item(27, 5) "  +4.63%  "

# Row 28
Set-TextValue $ws.Cells.Item(28, 4) "150.02"
Set-TextValue $ws.Cells.Item(28, 5) "  -1.47%  "

# Row 29
Set-TextValue $ws.Cells.Item(29, 4) "132.83"
Set-TextValue $ws.Cells.Item(29, 5) "  +3.92%  "

# Row 30
Set-TextValue $ws.Cells.Item(30, 2) "ImmutableX"
Set-TextValue $ws.Cells.Item(30, 3) "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Cells.Item(30, 4) "1.264"
Set-TextValue $ws.Cells.Item(30, 5) "  +29.85%  "

# Row 31
Set-TextValue $ws.Cells.Item(31, 2) "WrappedliquidstakedEther2.0"
Set-TextValue $ws.Cells.Item(31, 3) "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue $ws.Cells.Item(31, 4) "1.893.99"
Set-TextValue $ws.Cells.Item(31, 5) "  +1.65%  "

# Row 32
Set-TextValue $ws.Cells.Item(32, 4) "6.814"
Set-TextValue $ws.Cells.Item(32, 5) "  +8.38%  "

# Row 33
Set-TextValue $ws.Cells.Item(33, 4) "4.229"
Set-TextValue $ws.Cells.Item(33, 5) "  +4.46%  "

# Row 34
Set-TextValue $ws.Cells.Item(34, 5) "  +13.06%  "

# Row 35
Set-TextValue $ws.Cells.Item(35, 2) "WEMIXTOKEN"
Set-TextValue $ws.Cells.Item(35, 3) "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Cells.Item(35, 4) "1.780"
Set-TextValue $ws.Cells.Item(35, 5) "  +5.82%  "

# Row 36
Set-TextValue $ws.Cells.Item(36, 2) "Stellar"
Set-TextValue $ws.Cells.Item(36, 3) "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Cells.Item(36, 4) "0.08773"
Set-TextValue $ws.Cells.Item(36, 5) "  +4.05%  "

# Row 37
Set-TextValue $ws.Cells.Item(37, 4) "5.628"
Set-TextValue $ws.Cells.Item(37, 5) "  +5.67%  "

# Row 38
Set-TextValue $ws.Cells.Item(38, 4) "0.06665"
Set-TextValue $ws.Cells.Item(38, 5) "  +3.97%  "

# Row 39
Set-TextValue $ws.Cells.Item(39, 4) "9.165"
Set-TextValue $ws.Cells.Item(39, 5) "  +4.84%  "

# Row 40
Set-TextValue $ws.Cells.Item(40, 4) "0.02422"
Set-TextValue $ws.Cells.Item(40, 5) "  +4.33%  "

# Row 41
Set-TextValue $ws.Cells.Item(41, 4) "0.2252"
Set-TextValue $ws.Cells.Item(41, 5) "  +7.58%  "

# Row 42
Set-TextValue $ws.Cells.Item(42, 4) "1.270"
Set-TextValue $ws.Cells.Item(42, 5) "  +2.64%  "

# Row 43
Set-TextValue $ws.Cells.Item(43, 4) "0.6510"
Set-TextValue $ws.Cells.Item(43, 5) "  +6.60%  "

# Row 44
Set-TextValue $ws.Cells.Item(44, 4) "0.9996"
Set-TextValue $ws.Cells.Item(44, 5) "  +0.05%  "

# Row 45
Set-TextValue $ws.Cells.Item(45, 4) "13.85"
Set-TextValue $ws.Cells.Item(45, 5) "  +6.53%  "

# Row 46
Set-TextValue $ws.Cells.Item(46, 4) "0.6176"
Set-TextValue $ws.Cells.Item(46, 5) "  +4.74%  "

# Row 47
Set-TextValue $ws.Cells.Item(47, 4) "3.840"
Set-TextValue $ws.Cells.Item(47, 5) "  +2.47%  "

# Row 48
Set-TextValue $ws.Cells.Item(48, 4) "2.122"
Set-TextValue $ws.Cells.Item(48, 5) "  +5.10%  "

# Row 49
Set-TextValue $ws.Cells.Item(49, 4) "129.30"
Set-TextValue $ws.Cells.Item(49, 5) "  +2.08%  "

# Row 50
Set-TextValue $ws.Cells.Item(50, 4) "0.07327"
Set-TextValue $ws.Cells.Item(50, 5) "  +2.46%  "

# Row 51
Set-TextValue $ws.Cells.Item(51, 4) "80.41"
Set-TextValue $ws.Cells.Item(51, 5) "  +6.12%  "
